$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -0.2055599550297054
$ws.Range("C2").Value = -0.9145118796435601
$ws.Range("D2").Value = -0.422405887321083
$ws.Range("E2").Value = 0.04739885995190407
$ws.Range("F2").Value = 0.36022405491735
$ws.Range("G2").Value = -0.5181785263342299
$ws.Range("H2").Value = -1.248090541512699
$ws.Range("I2").Value = 0.4140885548769285
$ws.Range("J2").Value = -0.2051307335183153
$ws.Range("K2").Value = 0.03217303010139827

# Row 3
$ws.Range("B3").Value = -0.3319173426099191
$ws.Range("C3").Value = 0.160188649712558
$ws.Range("D3").Value = 0.6299933969855451
$ws.Range("E3").Value = 0.942818591950991
$ws.Range("F3").Value = 0.06441601069941108
$ws.Range("G3").Value = -0.6654960044790579
$ws.Range("H3").Value = 0.9966830919105695
$ws.Range("I3").Value = 0.3774638035153257
$ws.Range("J3").Value = 0.6147675671350392
$ws.Range("K3").Value = 0.5114901972596275

# Row 4
$ws.Range("B4").Value = 0.0645136089032002
$ws.Range("C4").Value = 0.5343183561761873
$ws.Range("D4").Value = 0.8471435511416332
$ws.Range("E4").Value = -0.03125903010994671
$ws.Range("F4").Value = -0.7611710452884157
$ws.Range("G4").Value = 0.9010080511012117
$ws.Range("H4").Value = 0.2817887627059679
$ws.Range("I4").Value = 0.5190925263256815
$ws.Range("J4").Value = 0.4158151564502698
$ws.Range("K4").Value = -0.3757832708791649

# Row 5
$ws.Range("B5").Value = 0.7947373931749101
$ws.Range("C5").Value = 1.107562588140356
$ws.Range("D5").Value = 0.2291600068887761
$ws.Range("E5").Value = -0.5007520082896928
$ws.Range("F5").Value = 1.161427088099934
$ws.Range("G5").Value = 0.5422077997046907
$ws.Range("H5").Value = 0.7795115633244043
$ws.Range("I5").Value = 0.6762341934489926
$ws.Range("J5").Value = -0.1153642338804421
$ws.Range("K5").Value = 0.4951246737870189

# Row 6
$ws.Range("B6").Value = 0.2721075919818648
$ws.Range("C6").Value = -0.6062949892697151
$ws.Range("D6").Value = -1.336207004448184
$ws.Range("E6").Value = 0.3259720919414433
$ws.Range("F6").Value = -0.2932471964538005
$ws.Range("G6").Value = -0.05594343283408693
$ws.Range("H6").Value = -0.1592208027094986
$ws.Range("I6").Value = -0.9508192300389333
$ws.Range("J6").Value = -0.3403303223714723
$ws.Range("K6").Value = -0.5018523531907899

# Row 7
$ws.Range("B7").Value = -0.5062156293670936
$ws.Range("C7").Value = -1.236127644545562
$ws.Range("D7").Value = 0.4260514518440648
$ws.Range("E7").Value = -0.193167836551179
$ws.Range("F7").Value = 0.04413592706853459
$ws.Range("G7").Value = -0.05914144280687711
$ws.Range("H7").Value = -0.8507398701363118
$ws.Range("I7").Value = -0.2402509624688508
$ws.Range("J7").Value = -0.4017729932881683
$ws.Range("K7").ClearContents()

# Row 8
$ws.Range("B8").Value = -0.8823411384658664
$ws.Range("C8").Value = 0.779837957923761
$ws.Range("D8").Value = 0.1606186695285172
$ws.Range("E8").Value = 0.3979224331482308
$ws.Range("F8").Value = 0.2946450632728191
$ws.Range("G8").Value = -0.4969533640566156
$ws.Range("H8").Value = 0.1135355436108454
$ws.Range("I8").Value = -0.04798648720847212
$ws.Range("J8").ClearContents()

# Row 9
$ws.Range("B9").Value = 0.6316605674913157
$ws.Range("C9").Value = 0.0124412790960719
$ws.Range("D9").Value = 0.2497450427157855
$ws.Range("E9").Value = 0.1464676728403738
$ws.Range("F9").Value = -0.6451307544890609
$ws.Range("G9").Value = -0.03464184682159993
$ws.Range("H9").Value = -0.1961638776409175
$ws.Range("I9").ClearContents()

# Row 10
$ws.Range("B10").Value = -0.1450047099080831
$ws.Range("C10").Value = 0.0922990537116305
$ws.Range("D10").Value = -0.0109783161637812
$ws.Range("E10").Value = -0.8025767434932158
$ws.Range("F10").Value = -0.1920878358257549
$ws.Range("G10").Value = -0.3536098666450724
$ws.Range("H10").ClearContents()

# Row 11
$ws.Range("B11").Value = 0.5929585102377013
$ws.Range("C11").Value = 0.4896811403622896
$ws.Range("D11").Value = -0.3019172869671451
$ws.Range("E11").Value = 0.3085716207003159
$ws.Range("F11").Value = 0.1470495898809984
$ws.Range("G11").ClearContents()

# Row 12
$ws.Range("B12").Value = 0.2093232598268204
$ws.Range("C12").Value = -0.5822751675026142
$ws.Range("D12").Value = 0.02821374016484672
$ws.Range("E12").Value = -0.1333082906544708
$ws.Range("F12").ClearContents()

# Row 13
$ws.Range("B13").Value = -0.4103003096576026
$ws.Range("C13").Value = 0.2001885980098584
$ws.Range("D13").Value = 0.03866656719054083
$ws.Range("E13").ClearContents()

# Row 14
$ws.Range("B14").Value = -0.1056739417364731
$ws.Range("C14").Value = -0.2671959725557906
$ws.Range("D14").ClearContents()

# Row 15
$ws.Range("B15").Value = 0.3451339801314955
$ws.Range("C15").ClearContents()

# Row 16
$ws.Range("B16").ClearContents()
